# ==========================================================================
# Edit script: Add Faculty details to Course_Information sheet and fix
# classroom allocations (C104 <-> C204 swap and related room corrections)
# across the workbook.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# 1. Course_Information: insert a "Faculty" column between
#    "Term Type"/"Basket" and "Display Format" (new column E, old E -> F)
# --------------------------------------------------------------------
$wsCourse = $wb.Worksheets.Item("Course_Information")

# Insert a new column at position 5 (E); existing column E (Display Format)
# shifts right to F automatically, along with the merged cells and the
# column F width (35) that was already set on old column E.
$wsCourse.Columns.Item(5).Insert()

# New column E should be 25 wide, keeping F at 35 (already carried over).
$wsCourse.Columns.Item(5).ColumnWidth = 25
$wsCourse.Columns.Item(6).ColumnWidth = 35

# Section headers for the new Faculty column
$wsCourse.Range("E4").Value = "Faculty"
$wsCourse.Range("E8").Value = "Faculty"

# Faculty values for the CORE COURSES section
$wsCourse.Range("E5").Value = "Animesh Roy, Pramod Yelmewad"

# Faculty values for the ELECTIVE COURSES section
$wsCourse.Range("E9").Value  = "Sandesh P"
$wsCourse.Range("E10").Value = "Shirshendu Layek"
$wsCourse.Range("E11").Value = "Krishnendu"
$wsCourse.Range("E12").Value = "Abdul Wahid"
$wsCourse.Range("E13").Value = "Malay Kumar"
$wsCourse.Range("E14").Value = "Sandesh Phalke"
$wsCourse.Range("E15").Value = "Anushree"
$wsCourse.Range("E16").Value = "Girish G N"
$wsCourse.Range("E17").Value = "Rajendra Hegadi"
$wsCourse.Range("E18").Value = "Sunil Saumya"
$wsCourse.Range("E19").Value = "Dibyajyothi"
$wsCourse.Range("E20").Value = "Chinmayananda A"
$wsCourse.Range("E21").Value = "Jagadish D N"
$wsCourse.Range("E22").Value = "Rajesh Kumar"
$wsCourse.Range("E23").Value = "Anand B"
$wsCourse.Range("E24").Value = "Aswath Babu"

# --------------------------------------------------------------------
# 2. Regular_Timetable / PreMid_Timetable / PostMid_Timetable:
#    room number corrections for the elective baskets
#    (same changes applied identically to all three sheets)
# --------------------------------------------------------------------
$timetableSheets = @("Regular_Timetable", "PreMid_Timetable", "PostMid_Timetable")
foreach ($sheetName in $timetableSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("D20").Value = "Mon 09:00-10:30 [C101], Wed 13:00-14:30 [C101]"
    $ws.Range("E20").Value = "Tue 14:30-15:30 [C101]"

    $ws.Range("D21").Value = "Mon 09:00-10:30 [C102], Wed 13:00-14:30 [C102]"
    $ws.Range("E21").Value = "Tue 14:30-15:30 [C102]"

    $ws.Range("D22").Value = "Mon 09:00-10:30 [C104], Wed 13:00-14:30 [C104]"
    $ws.Range("E22").Value = "Tue 14:30-15:30 [C104]"

    $ws.Range("D23").Value = "Tue 09:00-10:30 [C101], Thu 13:00-14:30 [C101]"
    $ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"

    $ws.Range("D24").Value = "Tue 09:00-10:30 [C102], Thu 13:00-14:30 [C102]"
    $ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"

    $ws.Range("D25").Value = "Tue 09:00-10:30 [C104], Thu 13:00-14:30 [C104]"
    $ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"

    $ws.Range("D26").Value = "Tue 09:00-10:30 [C202], Thu 13:00-14:30 [C202]"
    $ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"

    $ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"
    $ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"
    $ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"
    $ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"
    $ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"
}

# --------------------------------------------------------------------
# 3. Section_A: Mini Project room moved from C104 to C204
# --------------------------------------------------------------------
$wsSection = $wb.Worksheets.Item("Section_A")
$wsSection.Range("B16").Value = "Mini Project [C204]"
$wsSection.Range("C16").Value = "0-0-0-8-2 [C204]"
$wsSection.Range("D16").Value = "Full Sem [C204]"
$wsSection.Range("E16").Value = "0/0 [C204]"
$wsSection.Range("F16").Value = "0/0 [C204]"

# --------------------------------------------------------------------
# 4. Classroom_Utilization: swap weekly/daily hours between C104 and C204
# --------------------------------------------------------------------
$wsUtil = $wb.Worksheets.Item("Classroom_Utilization")
$wsUtil.Range("D9").Value = 0
$wsUtil.Range("E9").Value = 0
$wsUtil.Range("D16").Value = 7.5
$wsUtil.Range("E16").Value = 1.5

# --------------------------------------------------------------------
# 5. Classroom_Allocation: facility/room corrections
# --------------------------------------------------------------------
$wsAlloc = $wb.Worksheets.Item("Classroom_Allocation")

# Mini Project sessions move from C104 to C204, facility Projector -> TV
$wsAlloc.Range("I3").Value  = "TV"
$wsAlloc.Range("M3").Value  = "C204"
$wsAlloc.Range("I6").Value  = "TV"
$wsAlloc.Range("M6").Value  = "C204"
$wsAlloc.Range("I9").Value  = "TV"
$wsAlloc.Range("M9").Value  = "C204"
$wsAlloc.Range("I12").Value = "TV"
$wsAlloc.Range("M12").Value = "C204"
$wsAlloc.Range("I15").Value = "TV"
$wsAlloc.Range("M15").Value = "C204"

# CS468 Tue lecture room correction
$wsAlloc.Range("M31").Value = "C104"

# CS473 Tue lecture facility/room correction
$wsAlloc.Range("I32").Value = "Projector"
$wsAlloc.Range("M32").Value = "C202"

# DS456 Tue tutorial facility/room correction
$wsAlloc.Range("I33").Value = "Projector"
$wsAlloc.Range("M33").Value = "C101"

# EC456 Tue tutorial facility/room correction
$wsAlloc.Range("I34").Value = "Projector"
$wsAlloc.Range("M34").Value = "C102"

# DS401 Tue tutorial: room type/capacity/facility/room correction
$wsAlloc.Range("G35").Value = "classroom"
$wsAlloc.Range("H35").Value = "96"
$wsAlloc.Range("I35").Value = "Projector"
$wsAlloc.Range("M35").Value = "C104"

# PH454 Wed lecture facility/room correction
$wsAlloc.Range("I47").Value = "Projector"
$wsAlloc.Range("M47").Value = "C202"

# DE451 Wed lecture room correction
$wsAlloc.Range("M48").Value = "C203"

# DS456 Wed lecture room correction
$wsAlloc.Range("M49").Value = "C101"

# EC456 Wed lecture facility/room correction
$wsAlloc.Range("I50").Value = "Projector"
$wsAlloc.Range("M50").Value = "C102"

# DS401 Wed lecture facility/room correction
$wsAlloc.Range("I51").Value = "Projector"
$wsAlloc.Range("M51").Value = "C104"

# CS457 Wed tutorial facility/room correction
$wsAlloc.Range("I52").Value = "Projector"
$wsAlloc.Range("M52").Value = "C101"

# DS458 Wed tutorial facility/room correction
$wsAlloc.Range("I53").Value = "Projector"
$wsAlloc.Range("M53").Value = "C102"

# CS468 Wed tutorial facility/room correction
$wsAlloc.Range("I54").Value = "Projector"
$wsAlloc.Range("M54").Value = "C104"

# CS473 Wed tutorial room correction
$wsAlloc.Range("M55").Value = "C202"

# CS457 Thu lecture facility/room correction
$wsAlloc.Range("I60").Value = "Projector"
$wsAlloc.Range("M60").Value = "C101"

# DS458 Thu lecture facility/room correction
$wsAlloc.Range("I61").Value = "Projector"
$wsAlloc.Range("M61").Value = "C102"

# CS468 Thu lecture facility/room correction
$wsAlloc.Range("I62").Value = "Projector"
$wsAlloc.Range("M62").Value = "C104"

# CS473 Thu lecture room correction
$wsAlloc.Range("M63").Value = "C202"

# EC462 Thu tutorial facility/room correction
$wsAlloc.Range("I64").Value = "Projector"
$wsAlloc.Range("M64").Value = "C101"

# EC465 Thu tutorial room correction
$wsAlloc.Range("M65").Value = "C102"

# ASD352 Thu tutorial facility/room correction
$wsAlloc.Range("I66").Value = "Projector"
$wsAlloc.Range("M66").Value = "C104"

# PH454 Thu tutorial facility/room correction
$wsAlloc.Range("I67").Value = "Projector"
$wsAlloc.Range("M67").Value = "C202"

# DE451 Thu tutorial room correction
$wsAlloc.Range("M68").Value = "C203"

# --------------------------------------------------------------------
# 6. Basket_Course_Allocations: reduce allocated-rooms lists
# --------------------------------------------------------------------
$wsBasket = $wb.Worksheets.Item("Basket_Course_Allocations")
$wsBasket.Range("C2").Value  = "C004, C101"
$wsBasket.Range("C3").Value  = "C102"
$wsBasket.Range("C4").Value  = "C104"
$wsBasket.Range("C5").Value  = "C004, C101"
$wsBasket.Range("C6").Value  = "C102"
$wsBasket.Range("C7").Value  = "C104"
$wsBasket.Range("C8").Value  = "C202"
$wsBasket.Range("C9").Value  = "C203"
$wsBasket.Range("C10").Value = "C004, C101"
$wsBasket.Range("C11").Value = "C102"
$wsBasket.Range("C12").Value = "C104"
$wsBasket.Range("C13").Value = "C202"
